# Update "想去人数" (number of people interested) values for three events
# on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row 3 -> 南宁·草莓动漫节, row 5 -> 三月三国潮动漫节, row 6 -> 布谷鸟动漫展4th
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1871
$wsExhibit.Range("F5").Value = 671
$wsExhibit.Range("F6").Value = 231

# Sheet "全部类型": row 3 -> 南宁·草莓动漫节, row 6 -> 三月三国潮动漫节, row 7 -> 布谷鸟动漫展4th
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1871
$wsAll.Range("F6").Value = 671
$wsAll.Range("F7").Value = 231
